# Update countries & provincias Spain
# - Reorder "Bolivia" / "Republica de Yibuti" (Bolivia now ranks ahead) and
#   refresh their case counts.
# - Reorder "Guyana" into the Barbados/Mozambique block (Guyana now ranks
#   ahead of Mozambique/Haiti/San Martin) and refresh counts for that block.
# - Refresh Estados Unidos (row 4) and Brasil (row 14) case counts.
# - Bump the "Datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp (A1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 30 de Abril de 2020 a las 04:22"

# --- Estados Unidos (row 4) ---
$ws.Cells.Item(4, 2).Value = 1064533
$ws.Cells.Item(4, 3).Value = 339
$ws.Cells.Item(4, 5).Value = 855454
$ws.Cells.Item(4, 7).Value = 12
$ws.Cells.Item(4, 8).Value = 61668

# --- Brasil (row 14) ---
$ws.Cells.Item(14, 2).Value = 79685
$ws.Cells.Item(14, 3).Value = 324
$ws.Cells.Item(14, 5).Value = 40040
$ws.Cells.Item(14, 7).Value = 2
$ws.Cells.Item(14, 8).Value = 5513

# --- Bolivia / Republica de Yibuti swap (rows 88-89) ---
$ws.Cells.Item(88, 1).Value = "Bolivia"
$ws.Cells.Item(88, 2).Value = 1110
$ws.Cells.Item(88, 3).Value = 57
$ws.Cells.Item(88, 4).Value = 117
$ws.Cells.Item(88, 5).Value = 934
$ws.Cells.Item(88, 6).Value = 3
$ws.Cells.Item(88, 7).Value = 4
$ws.Cells.Item(88, 8).Value = 59

$ws.Cells.Item(89, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(89, 2).Value = 1077
$ws.Cells.Item(89, 4).Value = 599
$ws.Cells.Item(89, 5).Value = 476
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 8).Value = 2

# --- Guyana rotated ahead of Mozambique/Haiti/San Martin (rows 159-162) ---
$ws.Cells.Item(159, 1).Value = "Guyana"
$ws.Cells.Item(159, 2).Value = 78
$ws.Cells.Item(159, 3).Value = 4
$ws.Cells.Item(159, 4).Value = 18
$ws.Cells.Item(159, 5).Value = 52
$ws.Cells.Item(159, 6).Value = 5
$ws.Cells.Item(159, 8).Value = 8

$ws.Cells.Item(160, 1).Value = "Mozambique"
$ws.Cells.Item(160, 4).Value = 12
$ws.Cells.Item(160, 5).Value = 64
$ws.Cells.Item(160, 8).Value = 0

$ws.Cells.Item(161, 1).Value = "Haiti"
$ws.Cells.Item(161, 2).Value = 76
$ws.Cells.Item(161, 4).Value = 8
$ws.Cells.Item(161, 5).Value = 62
$ws.Cells.Item(161, 6).Value = 0
$ws.Cells.Item(161, 8).Value = 6

$ws.Cells.Item(162, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(162, 2).Value = 75
$ws.Cells.Item(162, 4).Value = 33
$ws.Cells.Item(162, 5).Value = 29
$ws.Cells.Item(162, 6).Value = 7
$ws.Cells.Item(162, 8).Value = 13
